# Applies the edit described by the commit:
# "organización de ekyword en test cases y steps, iteración con binding y cucumber"
#
# Summary of content changes on sheet1 ("Company Data"):
#  - Header row: column A changes from "XLS_PersonType" to "ID"
#    (the XLS_PersonType / Natural columns are dropped, replaced by a numeric ID column)
#  - Row 2: column A becomes the literal number 1 (was the text "Natural")
#    the rest of the row (B..J) keeps the same values
#  - A new Row 3 is added, identical to Row 2 except column A is the literal number 2,
#    including a mailto hyperlink on I3 (copied from I2) and the accompanying
#    relationship/hyperlink entry
#  - The used range / dimension grows from A1:J2 to A1:J3 and the active
#    selection moves to C12

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Company Data")

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "ID"

# --- Row 2: turn the person-type text into a numeric id ---------------
$ws.Range("A2").Value = 1

# --- Row 3: duplicate row 2, with id = 2 -------------------------------
$ws.Range("A2:J2").Copy($ws.Range("A3:J3"))
$ws.Range("A3").Value = 2

# Re-apply the hyperlink style + hyperlink relationship on I3 (mirrors I2)
$ws.Hyperlinks.Add($ws.Range("I3"), "mailto:mail@mail.com")
$ws.Range("I3").Style = "Hipervínculo"

# --- Selection ----------------------------------------------------------
$ws.Range("C12").Select() | Out-Null
